$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date the row was last changed.
# Update rows 2-20 from 2023-09-06 (45175) to 2023-09-08 (45177).
$ws.Range("C2:C20").Value = 45177
